# Apply edits to language.xlsx per commit:
# "stage 2 prep, some rearrangement of data, organism spawner, energy spawner within organism."
#
# Summary of change:
# - Remove standalone 'nucleoid' / 'ribosome' key rows (merge into 'essentialNucleoid' / 'essentialRibosome')
# - Remove 'bodyCorkscrew' / 'Corkscrew' pair, replace with 'bodySpirillum' / 'Spirillum'
# - Add new key/value pairs: essentialPlasmid / Plasmid DNA
# - Re-order many rows, add a new row 26
# - Update selection / view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New full data set (rows 2-26), columns A (Key) and B (Value).
# Rows 2-4 unchanged from before.
$data = @(
    @("welcome", "Welcome!"),
    @("title", "TBD"),
    @("none", "None"),
    @("test1", "Test 1"),
    @("test2", "Test 2"),
    @("testBodyCapsule", "Capsule"),
    @("testBodySphere", "Sphere"),
    @("categoryBody", "Body"),
    @("categoryCellStructure", "Cell Structure"),
    @("categoryMotility", "Motility"),
    @("essentialNucleoid", "Nucleoid"),
    @("essentialRibosome", "Ribosome"),
    @("essentialPlasmid", "Plasmid DNA"),
    @("bodyBacillus", "Bacillus"),
    @("bodyCoccus", "Coccus"),
    @("bodyCoccobacillus", "Coccobacillus"),
    @("bodySpirillum", "Spirillum"),
    @("cellStructureThermophile", "Thermophile"),
    @("cellStructurePsychrophile", "Psychrophile"),
    @("cellStructureMethanogen", "Methanogen"),
    @("cellStructureHalophile", "Halophile"),
    @("motilityFlagellaMonotrichous", "Monotrichous Flagella"),
    @("motilityFlagellaLophotrichous", "Lophotrichous Flagella"),
    @("motilityFlagellaPeritrichous", "Peritrichous Flagella"),
    @("motilityFlagellaAmphitrichous", "Amphitrichous Flagella")
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}

# Update view: clear topLeftCell scroll position and set new selection
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("A18").Select()
